# Adapt column header formatting to respective input file names (#7)
# - rename "*_old" headers to "*_FV2404"
# - rename "*_new" headers to "*_FV2410"
# - turn the A1:U57 range into an Excel Table ("Table1")
# - freeze the header row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header cells in row 1 (A1:U1)
for ($c = 1; $c -le 21; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $val = $cell.Text
    if ($val -ne $null) {
        if ($val.EndsWith("_old")) {
            $cell.Value = $val.Substring(0, $val.Length - 4) + "_FV2404"
        } elseif ($val.EndsWith("_new")) {
            $cell.Value = $val.Substring(0, $val.Length - 4) + "_FV2410"
        }
    }
}

# Convert range into a table
$tableRange = $ws.Range("A1:U57")
$table = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $tableRange, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$table.Name = "Table1"

# Freeze the header row (pane split after row 1)
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
